$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the underlying data values that drive the sheet's formulas and the
# two charts referencing Sheet1!$B$2:$B$21.
$ws.Range("B2").Value = 46
$ws.Range("B5").Value = 13
$ws.Range("B9").Value = 14
$ws.Range("B16").Value = 15
$ws.Range("B18").Value = 15

# Force a full recalculation so dependent formulas (C column percentages,
# B22/B24/B25 totals, C24/C25 sums) and the chart caches refresh.
$excel.CalculateFullRebuild()

# Restore the view: scroll back to the top-left and select B6 (matching the
# saved workbook view/selection state).
$ws.Activate()
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B6").Select()
